$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 101 (ALC)
$ws.Range("H101").Value = 203210.6
$ws.Range("I101").Value = 3622.6667
$ws.Range("J101").Value = 502592.5
$ws.Range("K101").Value = 10868.0001
$ws.Range("L101").Value = 1507777.5
$ws.Range("M101").Value = -9246.000100000001
$ws.Range("N101").Value = -1511021.5

# Row 113 (ALC)
$ws.Range("H113").Value = 3200.3845
$ws.Range("I113").Value = 3433.889
$ws.Range("J113").Value = 2675
$ws.Range("K113").Value = 3433.889
$ws.Range("L113").Value = 2675
$ws.Range("M113").Value = -179.8890000000001
$ws.Range("N113").Value = -9183

# Row 129 (ALC)
$ws.Range("H129").Value = 3704394.8
$ws.Range("I129").Value = 466
$ws.Range("J129").Value = 12346895
$ws.Range("K129").Value = 1398
$ws.Range("L129").Value = 37040685
$ws.Range("M129").Value = 3602
$ws.Range("N129").Value = -37050685

# Row 138 (ALC)
$ws.Range("H138").Value = 4463.61
$ws.Range("I138").Value = 1454.3572
$ws.Range("J138").Value = 4953.4883
$ws.Range("K138").Value = 4363.071599999999
$ws.Range("L138").Value = 14860.4649
$ws.Range("M138").Value = 776.9284000000007
$ws.Range("N138").Value = -25140.4649

$ws = $wb.Worksheets.Item("ARM")
# Row 124 (ARM)
$ws.Range("H124").Value = 18178.625
$ws.Range("J124").Value = 18178.625
$ws.Range("L124").Value = 18178.625
$ws.Range("N124").Value = -27998.625

# Row 125 (ARM)
$ws.Range("H125").Value = 32775.555
$ws.Range("J125").Value = 32775.555
$ws.Range("L125").Value = 32775.555
$ws.Range("N125").Value = -42615.555

$ws = $wb.Worksheets.Item("BSM")
# Row 27 (BSM)
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws.Range("H16").Value = 1233.3334
$ws.Range("J16").Value = 1250
$ws.Range("L16").Value = 1250
$ws.Range("N16").Value = -1824

# Row 31 (CRP)
$ws.Range("H31").Value = 2011.6316
$ws.Range("I31").Value = 1570.0952
$ws.Range("J31").Value = 2557.0588
$ws.Range("K31").Value = 1570.0952
$ws.Range("L31").Value = 2557.0588
$ws.Range("M31").Value = -1275.0952
$ws.Range("N31").Value = -3147.0588

# Row 34 (CRP)
$ws.Range("H34").Value = 2011.6316
$ws.Range("I34").Value = 1570.0952
$ws.Range("J34").Value = 2557.0588
$ws.Range("K34").Value = 1570.0952
$ws.Range("L34").Value = 2557.0588
$ws.Range("M34").Value = -1368.0952
$ws.Range("N34").Value = -2961.0588

# Row 86 (CRP)
$ws.Range("H86").Value = 3017.8333
$ws.Range("I86").Value = 3182.7727
$ws.Range("J86").Value = 2564.25
$ws.Range("K86").Value = 3182.7727
$ws.Range("L86").Value = 2564.25
$ws.Range("M86").Value = -2059.7727
$ws.Range("N86").Value = -4810.25

# Row 89 (CRP)
$ws.Range("H89").Value = 3017.8333
$ws.Range("I89").Value = 3182.7727
$ws.Range("J89").Value = 2564.25
$ws.Range("K89").Value = 15913.8635
$ws.Range("L89").Value = 12821.25
$ws.Range("M89").Value = -10297.8635
$ws.Range("N89").Value = -24053.25

# Row 94 (CRP)
$ws.Range("H94").Value = 1944.625
$ws.Range("I94").Value = 786.5714
$ws.Range("J94").Value = 2845.3333
$ws.Range("K94").Value = 786.5714
$ws.Range("L94").Value = 2845.3333
$ws.Range("M94").Value = -335.5714
$ws.Range("N94").Value = -3747.3333

# Row 105 (CRP)
$ws.Range("H105").Value = 2992
$ws.Range("I105").Value = 2253.3333
$ws.Range("K105").Value = 2253.3333
$ws.Range("M105").Value = -506.3332999999998

# Row 113 (CRP)
$ws.Range("H113").Value = 1233.3334
$ws.Range("J113").Value = 1250
$ws.Range("L113").Value = 1250
$ws.Range("N113").Value = -5590

# Row 140 (CRP)
$ws.Range("H140").Value = 55716
$ws.Range("J140").Value = 55716
$ws.Range("L140").Value = 55716
$ws.Range("N140").Value = -66076

$ws = $wb.Worksheets.Item("CUL")
# Row 68 (CUL)
$ws.Range("H68").Value = 1074.7128
$ws.Range("J68").Value = 1282.1875
$ws.Range("L68").Value = 3846.5625
$ws.Range("N68").Value = -5468.5625

# Row 71 (CUL)
$ws.Range("H71").Value = 1074.7128
$ws.Range("J71").Value = 1282.1875
$ws.Range("L71").Value = 11539.6875
$ws.Range("N71").Value = -19651.6875

# Row 107 (CUL)
$ws.Range("H107").Value = 883.4286
$ws.Range("I107").Value = 220.15384
$ws.Range("J107").Value = 1180.7587
$ws.Range("K107").Value = 660.4615200000001
$ws.Range("L107").Value = 3542.2761
$ws.Range("M107").Value = 1259.53848
$ws.Range("N107").Value = -7382.2761

# Row 113 (CUL)
$ws.Range("H113").Value = 102506.7
$ws.Range("I113").Value = 393.16666
$ws.Range("J113").Value = 116755.1
$ws.Range("K113").Value = 1179.49998
$ws.Range("L113").Value = 350265.3
$ws.Range("M113").Value = 990.5000199999999
$ws.Range("N113").Value = -354605.3

# Row 131 (CUL)
$ws.Range("H131").Value = 14556.507
$ws.Range("J131").Value = 1673.8871
$ws.Range("L131").Value = 5021.6613
$ws.Range("N131").Value = -15101.6613

$ws = $wb.Worksheets.Item("GSM")
# Row 123 (GSM)
$ws.Range("H123").Value = 18124.75
$ws.Range("J123").Value = 18124.75
$ws.Range("L123").Value = 18124.75
$ws.Range("N123").Value = -23024.75

# Row 138 (GSM)
$ws.Range("H138").Value = 19901
$ws.Range("J138").Value = 19901
$ws.Range("L138").Value = 19901
$ws.Range("N138").Value = -30181

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (LTW)
$ws.Range("H16").Value = 11675
$ws.Range("I16").Value = 16428.572
$ws.Range("J16").Value = 583.3333
$ws.Range("K16").Value = 16428.572
$ws.Range("L16").Value = 583.3333
$ws.Range("M16").Value = -16258.572
$ws.Range("N16").Value = -923.3333

# Row 93 (LTW)
$ws.Range("H93").Value = 1389.1333
$ws.Range("I93").Value = 1205.8889
$ws.Range("J93").Value = 1664
$ws.Range("K93").Value = 1205.8889
$ws.Range("L93").Value = 1664
$ws.Range("M93").Value = 42.11110000000008
$ws.Range("N93").Value = -4160

# Row 139 (LTW)
$ws.Range("H139").Value = 41479.445
$ws.Range("J139").Value = 41479.445
$ws.Range("L139").Value = 41479.445
$ws.Range("N139").Value = -51759.445

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (WVR)
$ws.Range("H62").Value = 3563.3333
$ws.Range("I62").Value = 3593.3333
$ws.Range("K62").Value = 3593.3333
$ws.Range("M62").Value = -2969.3333

# Row 65 (WVR)
$ws.Range("H65").Value = 3563.3333
$ws.Range("I65").Value = 3593.3333
$ws.Range("K65").Value = 17966.6665
$ws.Range("M65").Value = -14846.6665

# Row 74 (WVR)
$ws.Range("H74").Value = 10063
$ws.Range("J74").Value = 10063
$ws.Range("L74").Value = 10063
$ws.Range("N74").Value = -11935

# Row 77 (WVR)
$ws.Range("H77").Value = 10063
$ws.Range("J77").Value = 10063
$ws.Range("L77").Value = 30189
$ws.Range("N77").Value = -39549

# Row 81 (WVR)
$ws.Range("H81").Value = 1833.3334
$ws.Range("I81").Value = 1250
$ws.Range("K81").Value = 2500
$ws.Range("M81").Value = -1439

# Row 84 (WVR)
$ws.Range("H84").Value = 1833.3334
$ws.Range("I84").Value = 1250
$ws.Range("K84").Value = 12500
$ws.Range("M84").Value = -7196
